# Edit the "Web Parameters" sheet: adjust row 8/9/10/11 data and append
# two new question groups (rows 12-15), per commit:
#   "Added 2 more questions for calendar without interactoin."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web Parameters")

# ---- Row 8: amount_earlier 500 -> 300, amount_later 1000 -> 700 ----
$ws.Range("F8").Value = 300
$ws.Range("I8").Value = 700

# ---- Row 9: now the 2nd "calendar"/"none"/"none" sub-question ----
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "calendar"
$ws.Range("D9").Value = "none"
$ws.Range("E9").Value = "none"
$ws.Range("F9").Value = 500
$ws.Range("G9").ClearContents()
$ws.Range("H9").Value = 44593
$ws.Range("I9").Value = 800
$ws.Range("J9").ClearContents()
$ws.Range("K9").Value = 44617
$ws.Range("L9").Value = 1100
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = 100
$ws.Range("O9").Value = 100
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").Value = 8
$ws.Range("U9").Value = 8
$ws.Range("V9").Value = "Read 2001 example, absolute size"

# ---- Row 10: now the 3rd "calendar"/"none"/"none" sub-question ----
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "calendar"
$ws.Range("D10").Value = "none"
$ws.Range("E10").Value = "none"
$ws.Range("F10").Value = 300
$ws.Range("G10").ClearContents()
$ws.Range("H10").Value = 44593
$ws.Range("I10").Value = 1000
$ws.Range("J10").ClearContents()
$ws.Range("K10").Value = 44620
$ws.Range("L10").Value = 1100
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = 100
$ws.Range("O10").Value = 100
$ws.Range("T10").Value = 8
$ws.Range("U10").Value = 8
$ws.Range("V10").Value = "Read 2001 example, absolute size"

# ---- Row 11: treatment_id 6 -> 4, interaction titration -> drag ----
$ws.Range("A11").Value = 4
$ws.Range("D11").Value = "drag"

# ---- Row 12 (new): treatment_id 5, calendar/drag/laterAmount family ----
$ws.Range("A12").Value = 5
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "calendar"
$ws.Range("D12").Value = "drag"
$ws.Range("E12").Value = "laterAmount"
$ws.Range("F12").Value = 500
$ws.Range("H12").Value = 44593
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 44614
$ws.Range("L12").Value = 1100
$ws.Range("N12").Value = 100
$ws.Range("O12").Value = 100
$ws.Range("T12").Value = 8
$ws.Range("U12").Value = 8
$ws.Range("V12").Value = "Read 2001 example, absolute size"

# ---- Row 13 (new): treatment_id 6, word/titration/laterAmount ----
$ws.Range("A13").Value = 6
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "word"
$ws.Range("D13").Value = "titration"
$ws.Range("E13").Value = "laterAmount"
$ws.Range("F13").Value = 500
$ws.Range("G13").Value = 2
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 3
$ws.Range("M13").Value = 10
$ws.Range("V13").Value = "Read 2001 example, absolute size"

# ---- Row 14 (new): treatment_id 7, barchart/titration/laterAmount ----
$ws.Range("A14").Value = 7
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "barchart"
$ws.Range("D14").Value = "titration"
$ws.Range("E14").Value = "laterAmount"
$ws.Range("F14").Value = 500
$ws.Range("G14").Value = 2
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 10
$ws.Range("L14").Value = 1500
$ws.Range("M14").Value = 10
$ws.Range("N14").Value = 100
$ws.Range("O14").Value = 100
$ws.Range("P14").Value = 0.5
$ws.Range("Q14").Value = 0.5
$ws.Range("R14").Value = 8
$ws.Range("S14").Value = 8
$ws.Range("T14").Value = 8.5
$ws.Range("U14").Value = 8.5
$ws.Range("V14").Value = "Read 2001 example, absolute size"

# ---- Row 15 (new): treatment_id 8, calendar/titration/laterAmount ----
$ws.Range("A15").Value = 8
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "calendar"
$ws.Range("D15").Value = "titration"
$ws.Range("E15").Value = "laterAmount"
$ws.Range("F15").Value = 500
$ws.Range("H15").Value = 44593
$ws.Range("I15").Value = 1000
$ws.Range("K15").Value = 44614
$ws.Range("L15").Value = 1100
$ws.Range("N15").Value = 100
$ws.Range("O15").Value = 100
$ws.Range("T15").Value = 8
$ws.Range("U15").Value = 8
$ws.Range("V15").Value = "Read 2001 example, absolute size"

# Re-apply the existing date format (copied from H8/K8, which already use
# the workbook's date style) to the new/re-purposed date_earlier /
# date_later cells, so they share the same style index instead of Excel
# minting a brand-new numFmt. (Multi-area ranges only paste into the first
# area here, so paste one cell at a time.)
foreach ($cellRef in @("H9", "H10", "H12", "H15")) {
    $ws.Range("H8").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}
foreach ($cellRef in @("K9", "K10", "K12", "K15")) {
    $ws.Range("K8").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Move selection to match the post-edit state (P10 was last active cell).
$ws.Range("P10").Select()
